$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L3").Value = 1.32
$ws.Range("F4").Value = 1.8
$ws.Range("G4").Value = 2.04
$ws.Range("H4").Value = 4
$ws.Range("O4").Value = 1.26
$ws.Range("W4").Value = 1.97
$ws.Range("J6").Value = 6.2
$ws.Range("N7").Value = 1.36
$ws.Range("P7").Value = 1.36
$ws.Range("T7").Value = 1.04
$ws.Range("U7").Value = 1.04
$ws.Range("L8").Value = 1.45
$ws.Range("N8").Value = 3.5
$ws.Range("O8").Value = 1.38
$ws.Range("R8").Value = 1.31
$ws.Range("V8").Value = 1.37
$ws.Range("W8").Value = 1.76
$ws.Range("Z8").Value = 24
$ws.Range("AB8").Value = 9
$ws.Range("AC8").Value = 7.4
$ws.Range("AF8").Value = 13
$ws.Range("AG8").Value = 11
$ws.Range("AI8").Value = 60
$ws.Range("AJ8").Value = 29
$ws.Range("AK8").Value = 25
$ws.Range("AM8").Value = 120
$ws.Range("AO8").Value = 48
$ws.Range("L9").Value = 1.32
$ws.Range("S9").Value = 2.74
$ws.Range("V9").Value = 1.17
$ws.Range("W9").Value = 2.72
$ws.Range("X9").Value = 21
$ws.Range("Y9").Value = 26
$ws.Range("AG9").Value = 9.4
$ws.Range("AJ9").Value = 14.5
$ws.Range("AM9").Value = 85
$ws.Range("G10").Value = 2.1
$ws.Range("H10").Value = 3.5
$ws.Range("I10").Value = 4.6
$ws.Range("L10").Value = 1.01
$ws.Range("M10").Value = 1.01
$ws.Range("N10").Value = 2.5
$ws.Range("O10").Value = 1.17
$ws.Range("R10").Value = 1.52
$ws.Range("S10").Value = 2.1
$ws.Range("T10").Value = 1.01
$ws.Range("U10").Value = 1.01
$ws.Range("V10").Value = 1.27
$ws.Range("W10").Value = 1.9
$ws.Range("X10").Value = 1000
$ws.Range("Y10").Value = 1000
$ws.Range("Z10").Value = 1000
$ws.Range("AA10").Value = 1000
$ws.Range("AB10").Value = 1000
$ws.Range("AC10").Value = 1000
$ws.Range("AD10").Value = 1000
$ws.Range("AE10").Value = 1000
$ws.Range("AF10").Value = 1000
$ws.Range("AG10").Value = 1000
$ws.Range("AH10").Value = 1000
$ws.Range("AI10").Value = 1000
$ws.Range("AJ10").Value = 1000
$ws.Range("AK10").Value = 1000
$ws.Range("AL10").Value = 1000
$ws.Range("AM10").Value = 1000
$ws.Range("AN10").Value = 1000
$ws.Range("AO10").Value = 1000
$ws.Range("N11").Value = 3.8
$ws.Range("P11").Value = 1.93
$ws.Range("AI11").Value = 60
$ws.Range("AK11").Value = 26
$ws.Range("F12").Value = 8
$ws.Range("S12").Value = 2.78
$ws.Range("T12").Value = 1.92
$ws.Range("AC12").Value = 11
